# Split "Ministry Course Code and Level" column into two columns:
# "Ministry Course Code" and "Ministry Course Level" (Summer Reporting File Spec)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "Session Date" column (H), which pushes
# Session Date / Final Percent / Final Letter Grade / Credits one column to
# the right and gives us a second column to hold the split-out course level.
$ws.Columns("H").Insert()

# --- Header + data: split "ENST 12" into code "ENST" and level 12 ---
$ws.Range("G1").Value = "Ministry Course Code"

$ws.Range("G2").Value = "ENST"
$ws.Range("G3").Value = "ENST"
$ws.Range("G4").Value = "ENST"

$ws.Range("H1").Value = "Ministry Course Level"

$ws.Range("H2").Value = 12
$ws.Range("H3").Value = 12
$ws.Range("H4").Value = 12

# --- Sheet formatting tweaks ---
$ws.Rows(1).RowHeight = 43.2
$ws.Columns("H").ColumnWidth = 12

# --- Selection matches the target worksheet view ---
$ws.Range("G1:H1048576").Select()
